$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sellers = @{
    2  = "Ellie Ellie"
    3  = "thechucklingcheesecompany"
    4  = "Dibor"
    5  = "Gaamaa"
    6  = "Oakdene Designs"
    7  = "My Posh Shop"
    8  = "lovetreedesign"
    9  = "madewithlovecardboutique"
    10 = "songsofinkandsteel"
    11 = "thegourmetchocolatepizzaco"
    12 = "qwertybeerbox"
    13 = "The Alphabet Gift Shop"
    14 = "ladedaliving"
    15 = "Hurleyburley man"
    16 = "Joy by Corrine Smith"
    17 = "The Rustic Dish®"
    18 = "Lisa Angel"
    19 = "The Forest & Co"
    20 = "Hurleyburley man"
    21 = "alphabetinteriors"
    22 = "Dibor"
}

foreach ($row in $sellers.Keys) {
    $ws.Range("C$row").Value = $sellers[$row]
}
